$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (trial numbers) for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) updated values for columns B:E
$ws.Range("B2").Value = 12.232762440473394
$ws.Range("C2").Value = 11.559030914856949
$ws.Range("D2").Value = 13.066501868432459
$ws.Range("E2").Value = 12.322849670143519

# Row 3 (STR) updated values for columns B:E
$ws.Range("B3").Value = 11.090405638127915
$ws.Range("C3").Value = 10.465699751611059
$ws.Range("D3").Value = 12.552884336383746
$ws.Range("E3").Value = 11.543594037139135

# Update selection to match the new, narrower active region
$ws.Range("B1:E3").Select()
